$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(18).Insert()

$ws.Range("A18").Value = 5
$ws.Range("B18").Value = "O equipamento agora está funcionando perfeitamente. Ainda não testamos para a esterilização do nosso produto. Mas acredito que vai da tudo certo."
$ws.Range("C18").Value = 45954.49367453704
$ws.Range("D18").Value = "ZjIwMzgwNmItMDdmNy00OWM3LWIxMzEtYTVhZjkwNDM0OGQzOjU3MDE2"
